$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 0. Capture a style-only donor for style index 6 (only currently present on row 12) ---
# before row 12 gets overwritten, copy its formatting to a scratch row far away.
$ws.Range("A12:D12").Copy()
$ws.Range("A40").PasteSpecial(-4122)   # xlPasteFormats

# --- 1. New row 15 := old row 11 ("Note" comment cell), moved down ---
$ws.Range("A11").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats (style s=1)
$ws.Range("A15").Value = "Note"

# --- 2. Row 11 becomes a brand new row: only B11 has content, A11 is wiped ---
$ws.Range("A11").Clear()
$ws.Range("B11").Value = "Line with null cell"

# --- 3. New row 12 := old row 13 (style s=2, all empty) ---
$ws.Range("A13:D13").Copy()
$ws.Range("A12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A12:D12").ClearContents()

# --- 4. New row 13 := old row 14 (style s=4, all empty) ---
$ws.Range("A14:D14").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A13:D13").ClearContents()

# --- 5. New row 14 := old row 12 (style s=6, values 78174 / "String" / formula) ---
$ws.Range("A40:D40").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A14:D14").ClearContents()
$ws.Range("A14").Value = 78174
$ws.Range("B14").Value = "String"
$ws.Range("D14").Formula = "=A2+A2"

# --- 6. New row 16: B16 new string, C16 empty styled cell (s=2) ---
$ws.Range("B16").Value = "Line with blank after last non-empty"
$ws.Range("A1").Copy()
$ws.Range("C16").PasteSpecial(-4122)   # xlPasteFormats (style s=2)
$ws.Range("C16").ClearContents()

# --- cleanup scratch area ---
$ws.Range("A40:D40").Clear()

# --- comment: move from A11 to A15 (the A11 comment was already dropped by Clear() above) ---
$ws.Range("A15").AddComment("Uma nota de exemplo")

# --- column B width (closest achievable to 33.85546875 chars) ---
$ws.Columns("B:B").ColumnWidth = 33

# --- selection ---
$ws.Range("C16").Select()
